# Add the "2022-Q4" data: a new worksheet with fund-holding detail, plus a
# new summary row on the "总计" sheet. The new sheet is inserted right
# after "总计" and before the existing "2022-Q2" sheet; "2022-Q2" and
# "2020-Q4" simply shift one tab position to the right.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the range to keep its literal text (otherwise Excel would
    # auto-convert a numeric-looking string like "009686" to the number
    # 9686 and drop the leading zero).
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert a new row for 2022-Q4 above the existing rows.
#    Existing rows 2 ("2022-Q2") and 3 ("2020-Q4") shift down to rows 3
#    and 4 respectively, and the running index in column A is renumbered.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Copy column-A's numbered-index style down onto the new row 4 first, so
# the appended row matches the existing formatting (bold/centered/bordered).
$summary.Range("A3").Copy()
$summary.Range("A4").PasteSpecial(-4122)

# Write bottom-up so we never clobber a value before reading it.
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2020-Q4"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.01

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 3
$summary.Range("D3").Value = 0.01

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.65

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计" (i.e. before
#    the current "2022-Q2" sheet), and populate it with fund data.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Reuse the existing bold/bordered styles (from "总计") for the header
# row and the running-index column, instead of inventing new ones.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$summary.Range("A2").Copy()
$q4.Range("A2:A5").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Row 2
$q4.Range("A2").Value = 0
Set-TextValue $q4.Range("B2") "009686"
Set-TextValue $q4.Range("C2") "华夏磐利一年定期开放混合A"
Set-TextValue $q4.Range("D2") "11.49"
Set-TextValue $q4.Range("E2") "92.56"
Set-TextValue $q4.Range("F2") "4.23"
Set-TextValue $q4.Range("G2") "0.4860"
$q4.Range("H2").Value = 5

# Row 3
$q4.Range("A3").Value = 1
Set-TextValue $q4.Range("B3") "015697"
Set-TextValue $q4.Range("C3") "华夏磐润两年定开混合A"
Set-TextValue $q4.Range("D3") "2.68"
Set-TextValue $q4.Range("E3") "86.76"
Set-TextValue $q4.Range("F3") "4.07"
Set-TextValue $q4.Range("G3") "0.1091"
$q4.Range("H3").Value = 7

# Row 4
$q4.Range("A4").Value = 2
Set-TextValue $q4.Range("B4") "015698"
Set-TextValue $q4.Range("C4") "华夏磐润两年定开混合C"
Set-TextValue $q4.Range("D4") "0.99"
Set-TextValue $q4.Range("E4") "86.76"
Set-TextValue $q4.Range("F4") "4.07"
Set-TextValue $q4.Range("G4") "0.0403"
$q4.Range("H4").Value = 7

# Row 5
$q4.Range("A5").Value = 3
Set-TextValue $q4.Range("B5") "009687"
Set-TextValue $q4.Range("C5") "华夏磐利一年定期开放混合C"
Set-TextValue $q4.Range("D5") "0.46"
Set-TextValue $q4.Range("E5") "92.56"
Set-TextValue $q4.Range("F5") "4.23"
Set-TextValue $q4.Range("G5") "0.0195"
$q4.Range("H5").Value = 5

# Keep the original active/selected tab ("2020-Q4") instead of leaving the
# newly-inserted sheet selected.
$wb.Worksheets.Item("2020-Q4").Activate()
